# Rename worksheets in-place (positional 1:1 mapping) to new summ<id> names.
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ59379425",
    "summ59527898",
    "summ59712032",
    "summ59912815",
    "summ00180925",
    "summ00539837",
    "summ00776111",
    "summ00976166",
    "summ01178229",
    "summ01365793",
    "summ01607572",
    "summ01875758",
    "summ02087877",
    "summ02287920",
    "summ02475968",
    "summ02668004",
    "summ02856181",
    "summ03123161",
    "summ03298840",
    "summ03494011",
    "summ03696043",
    "summ03886316",
    "summ04078345",
    "summ04260115",
    "summ04479019",
    "summ04690057",
    "summ04901612",
    "summ05097638",
    "summ05297840",
    "summ05506937",
    "summ05696965",
    "summ05894594",
    "summ06090116",
    "summ06292465",
    "summ06487496",
    "summ06754531",
    "summ06956071",
    "summ07195191",
    "summ07392596",
    "summ07593628",
    "summ07779998",
    "summ07968951",
    "summ08211110",
    "summ08399779",
    "summ08580801",
    "summ08766827",
    "summ08946663",
    "summ09145143",
    "summ09331177",
    "summ09525214"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
